$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update MTC column for rows 8-12 (new MTC test case IDs)
$ws.Range("B8").Value = "MTC_007"
$ws.Range("B9").Value = "MTC_008"
$ws.Range("B10").Value = "MTC_009"
$ws.Range("B11").Value = "MTC_010"
$ws.Range("B12").Value = "MTC_011"

# Update Execute column for rows 9 and 11 to "Yes"
$ws.Range("E9").Value = "Yes"
$ws.Range("E11").Value = "Yes"

# Reset the custom font formatting on E1 (was JetBrains Mono) back to the
# same style used by the rest of the header row
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Update the active selection to B9:B12
$ws.Range("B9:B12").Select()
